# Insert a new weekly price record as row 139, shifting the existing
# rows 139-202 down to 140-203 (the newest reading moves to the top of
# this data block, pushing the rest of the series down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("139").Insert()

$ws.Range("A139").Value = 5
$ws.Range("B139").Value = "Macroferia Regional de Talca"
$ws.Range("C139").Value = "Maule"
$ws.Range("D139").Value = 44510
$ws.Range("E139").Value = 7
$ws.Range("F139").Value = 100114014
$ws.Range("G139").Value = "Betarraga"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 6000
$ws.Range("K139").Value = 550
$ws.Range("L139").Value = 550
$ws.Range("M139").Value = 550
$ws.Range("N139").Value = "`$/paquete 5 unidades"
$ws.Range("O139").Value = "Región del Maule"
$ws.Range("P139").Value = 110
$ws.Range("Q139").Value = 5
$ws.Range("R139").Value = "Hortaliza"
